# Apply "Initial results in tex." edit to the Summary sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$cities = @('Bath','Belfast','Brighton','Bristol','Cardiff','Coventry','Exeter','Glasgow','Leeds','Leicester',
            'Liverpool','Manchester','Newcastle','Nottingham','Oxford','Plymouth','Sheffield','Southampton',
            'Sunderland','York','Belgrade','Berlin','Boston','Dublin','Minsk')

$header = "City & Best & Avg. & T & Alg. & Best & Avg. \\ \hline"

$blocks = @(4, 34, 64)

# First: rewrite column A with the plain city names (this introduces the new
# shared strings "Bath".."Minsk" in that order, matching the first block).
foreach ($blockStart in $blocks) {
    for ($i = 0; $i -lt $cities.Length; $i++) {
        $row = $blockStart + $i
        $ws.Range("A$row").Value = $cities[$i]
    }
}

# Second: the "\hline" totals-row marker (new shared string right after the
# city names).
foreach ($blockStart in $blocks) {
    $totalsRow = $blockStart + 25
    $ws.Range("P$totalsRow").Value = "\hline"
}

# Third: the LaTeX table header text (last new shared string).
foreach ($blockStart in $blocks) {
    $headerRow = $blockStart - 1
    $ws.Range("P$headerRow").Value = $header
}

# Finally: the per-row LaTeX formulas in column P.
foreach ($blockStart in $blocks) {
    for ($i = 0; $i -lt $cities.Length; $i++) {
        $row = $blockStart + $i
        $formula = '=A' + $row + '&"&"&IF(E' + $row + '<=L' + $row + ',"\bf{"&E' + $row + '&"}",E' + $row + ')&"&"&ROUND(G' + $row + ',1)&"&"&ROUND(J' + $row + ',1)&"&"&K' + $row + '&"&"&IF(L' + $row + '<=E' + $row + ',"\bf{"&L' + $row + '&"}",L' + $row + ')&"&"&ROUND(M' + $row + ',1)&"\\"'
        $ws.Range("P$row").Formula = $formula
    }
}

# Update the sheet view to match the committed state.
$ws.Application.ActiveWindow.Zoom = 115
$ws.Application.ActiveWindow.ScrollRow = 66
$ws.Range("P64:P88").Select()
